$wb = $excel.ActiveWorkbook

# --- 1. "Units" sheet: record vom_cost_Input1 = 1 for the Methanol_Reactor row ---
$wsUnits = $wb.Worksheets.Item("Units")
$wsUnits.Range("U3").Value = 1

# --- 2. "Connections" sheet: the per-input/output vom_cost split replaces the
#        single combined "vom_cost" column, so drop that column from the table. ---
$wsConn = $wb.Worksheets.Item("Connections")
$table = $wsConn.ListObjects.Item("Table13")

# Removing the whole sheet column shifts everything after it (and the data
# cells) one slot to the left, which is what turns the old
# vom_cost / vom_cost_Input1 / vom_cost_Input2 / vom_cost_Output1 / vom_cost_Output2
# run into vom_cost_Input1 / vom_cost_Input2 / vom_cost_Output1 / vom_cost_Output2.
$wsConn.Range("V1:V5").EntireColumn.Delete()

# Shrink the table definition to match the new extent.
$table.Resize($wsConn.Range("A1:Y5"))

# Re-assert each shifted header's value so the table's column-name metadata
# (which isn't touched by a plain column delete) re-syncs with the sheet.
foreach ($colLetter in @("V", "W", "X", "Y")) {
    $headerCell = $wsConn.Range($colLetter + "1")
    $headerCell.Value = $headerCell.Value()
}

# --- 3. Restore the selections left behind by the edit ---
$wsUnits.Range("R23").Select()
$wsConn.Range("V3").Select()
$wsUnits.Activate()
